$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sign-in/out log entries appended below the existing data (rows 4-7).
# Columns: A=Student Number, B=First Name, C=Last Name, D=Date,
#          E=Sign-In Time, F=Sign-Out Time, G=Teacher, H=Reason
#
# Student Number and Date are stored as plain text in this sheet (matching
# rows 2-3), so the all-digit / date-like strings below must be forced to
# Text format first -- otherwise Excel would auto-convert them to a number
# or a date serial value. The format is reset back to Normal afterwards so
# the cells end up styled the same as the rest of the sheet.
$numberCol = $ws.Range("A4:A7")
$dateCol = $ws.Range("D4:D7")
$numberCol.NumberFormat = "@"
$dateCol.NumberFormat = "@"

$ws.Range("A4").Value = "111111111"
$ws.Range("B4").Value = "Katelyn"
$ws.Range("C4").Value = "W"
$ws.Range("D4").Value = "2018/11/18"
$ws.Range("E4").Value = "4:54 PM"
$ws.Range("G4").Value = "Math"
$ws.Range("H4").Value = "Quiet Work"

$ws.Range("A5").Value = "111111111"
$ws.Range("B5").Value = "Yash"
$ws.Range("C5").Value = "A"
$ws.Range("D5").Value = "2018/11/18"
$ws.Range("E5").Value = "4:58 PM"
$ws.Range("G5").Value = "Music"
$ws.Range("H5").Value = "Quiet Work"

$ws.Range("A6").Value = "222222222"
$ws.Range("B6").Value = "Guy"
$ws.Range("C6").Value = "M"
$ws.Range("D6").Value = "2018/11/18"
$ws.Range("E6").Value = "4:58 PM"
$ws.Range("G6").Value = "Math"
$ws.Range("H6").Value = "Academic Support"

$ws.Range("A7").Value = "111111111"
$ws.Range("B7").Value = "Yash"
$ws.Range("C7").Value = "A"
$ws.Range("D7").Value = "2018/11/18"
$ws.Range("E7").Value = "4:59 PM"
$ws.Range("G7").Value = "Math"
$ws.Range("H7").Value = "Chill Zone"

$numberCol.Style = "Normal"
$dateCol.Style = "Normal"
